# PBAC WIP Phase 2
# Insert a new "Fedlead-positive" column before the existing
# "Permission Internal Name" column (old column N becomes O).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at N; this shifts the old N ("Permission Internal
# Name") to O, along with its formatting/width.
$ws.Columns("N").Insert()

# Header for the newly inserted column.
$ws.Range("N1").Value = "Fedlead-positive"

# Column N width (old N keeps its width, which moved to O automatically).
$ws.Columns("N").ColumnWidth = 14

# Per-row values for the new column.
$values = @{
    2  = "checked"
    3  = "checked"
    4  = "checked"
    5  = "checked"
    6  = "checked"
    7  = "checked"
    8  = "checked"
    9  = "checked"
    10 = "checked"
    11 = "fixed_unchecked"
    12 = "checked"
    13 = "checked"
    14 = "checked"
    15 = "checked"
    16 = "checked"
    17 = "checked"
    18 = "fixed_unchecked"
    19 = "checked"
    20 = "checked"
    21 = "checked"
    22 = "checked"
    23 = "checked"
    24 = "checked"
    25 = "checked"
    26 = "checked"
    27 = "checked"
    28 = "checked"
    29 = "checked"
    30 = "checked"
    31 = "checked"
    32 = "checked"
    33 = "checked"
    34 = "fixed_unchecked"
    35 = "checked"
    36 = "fixed_checked"
    37 = "fixed_unchecked"
}

foreach ($row in $values.Keys) {
    $ws.Range("N$row").Value = $values[$row]
}

# Update the active selection to match the saved workbook state.
$ws.Range("N42").Select()
